$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2108980827447023
$ws.Range("C2").Value = 0.515640766902119
$ws.Range("J2").Value = 0.01614530776992936
$ws.Range("P2").Value = 0.1493440968718466
$ws.Range("S2").Value = 0.1079717457114026
$ws.Range("B3").Value = 0.01145038167938931
$ws.Range("C3").Value = 0.03244274809160305
$ws.Range("J3").Value = 0.0267175572519084
$ws.Range("P3").Value = 0.6984732824427481
$ws.Range("S3").Value = 0.2309160305343511
$ws.Range("J4").Value = 0.0975609756097561
$ws.Range("P4").Value = 0.6829268292682927
$ws.Range("S4").Value = 0.2195121951219512
$ws.Range("O5").Value = 0.1
$ws.Range("P5").Value = 0.6
$ws.Range("S5").Value = 0.3
$ws.Range("B6").Value = 0.07272727272727272
$ws.Range("D6").Value = 0.01363636363636364
$ws.Range("E6").Value = 0.001515151515151515
$ws.Range("F6").Value = 0.06363636363636363
$ws.Range("J6").Value = 0.2363636363636364
$ws.Range("O6").Value = 0.0303030303030303
$ws.Range("Q6").Value = 0.1287878787878788
$ws.Range("R6").Value = 0.07272727272727272
$ws.Range("S6").Value = 0.3803030303030303
$ws.Range("B7").Value = 0.088379705400982
$ws.Range("D7").Value = 0.01800327332242226
$ws.Range("F7").Value = 0.04746317512274959
$ws.Range("J7").Value = 0.1358428805237316
$ws.Range("O7").Value = 0.02127659574468085
$ws.Range("Q7").Value = 0.1636661211129296
$ws.Range("R7").Value = 0.07037643207855974
$ws.Range("S7").Value = 0.4549918166939443
$ws.Range("B8").Value = 0.1002017484868863
$ws.Range("D8").Value = 0.01412239408204439
$ws.Range("E8").Value = 0.0006724949562878278
$ws.Range("F8").Value = 0.05110961667787491
$ws.Range("J8").Value = 0.1149966375252186
$ws.Range("O8").Value = 0.02017484868863484
$ws.Range("Q8").Value = 0.1687962340282448
$ws.Range("R8").Value = 0.09952925353059852
$ws.Range("S8").Value = 0.4303967720242098
$ws.Range("B9").Value = 0.09171075837742504
$ws.Range("D9").Value = 0.01763668430335097
$ws.Range("E9").Value = 0.003527336860670194
$ws.Range("F9").Value = 0.06701940035273368
$ws.Range("J9").Value = 0.1234567901234568
$ws.Range("O9").Value = 0.01763668430335097
$ws.Range("Q9").Value = 0.1569664902998236
$ws.Range("R9").Value = 0.07936507936507936
$ws.Range("S9").Value = 0.4426807760141093
$ws.Range("B10").Value = 0.1131520940484938
$ws.Range("D10").Value = 0.01959343619887338
$ws.Range("E10").Value = 0.001959343619887338
$ws.Range("F10").Value = 0.06955669850600049
$ws.Range("J10").Value = 0.1234386480529023
$ws.Range("O10").Value = 0.01542983100661278
$ws.Range("Q10").Value = 0.2192015674748959
$ws.Range("R10").Value = 0.07323046779328925
$ws.Range("S10").Value = 0.3644379132990448
$ws.Range("F11").Value = 0.0009775171065493646
$ws.Range("G11").Value = 0.155425219941349
$ws.Range("J11").Value = 0.09286412512218964
$ws.Range("K11").Value = 0.2160312805474096
$ws.Range("L11").Value = 0.5190615835777126
$ws.Range("S11").Value = 0.01564027370478983
$ws.Range("G12").Value = 0.7148014440433214
$ws.Range("J12").Value = 0.2256317689530686
$ws.Range("K12").Value = 0.007220216606498195
$ws.Range("L12").Value = 0.02166064981949458
$ws.Range("S12").Value = 0.03068592057761733
$ws.Range("F13").Value = 0.006993006993006993
$ws.Range("G13").Value = 0.5664335664335665
$ws.Range("J13").Value = 0.3566433566433567
$ws.Range("S13").Value = 0.06993006993006994
$ws.Range("F15").Value = 0.01481481481481482
$ws.Range("H15").Value = 0.1496296296296296
$ws.Range("I15").Value = 0.05037037037037037
$ws.Range("J15").Value = 0.362962962962963
$ws.Range("K15").Value = 0.06518518518518518
$ws.Range("M15").Value = 0.01037037037037037
$ws.Range("O15").Value = 0.05333333333333334
$ws.Range("S15").Value = 0.2933333333333333
$ws.Range("F16").Value = 0.01541095890410959
$ws.Range("H16").Value = 0.1832191780821918
$ws.Range("I16").Value = 0.0821917808219178
$ws.Range("J16").Value = 0.3938356164383562
$ws.Range("K16").Value = 0.1027397260273973
$ws.Range("M16").Value = 0.02397260273972603
$ws.Range("O16").Value = 0.0410958904109589
$ws.Range("S16").Value = 0.1575342465753425
$ws.Range("F17").Value = 0.01343705799151344
$ws.Range("H17").Value = 0.1612446958981612
$ws.Range("I17").Value = 0.08345120226308345
$ws.Range("J17").Value = 0.4243281471004243
$ws.Range("K17").Value = 0.1025459688826025
$ws.Range("M17").Value = 0.01626591230551627
$ws.Range("O17").Value = 0.05657708628005657
$ws.Range("S17").Value = 0.1421499292786421
$ws.Range("F18").Value = 0.02072538860103627
$ws.Range("H18").Value = 0.155440414507772
$ws.Range("I18").Value = 0.08290155440414508
$ws.Range("J18").Value = 0.4127806563039724
$ws.Range("K18").Value = 0.1174438687392055
$ws.Range("M18").Value = 0.01381692573402418
$ws.Range("O18").Value = 0.0535405872193437
$ws.Range("S18").Value = 0.1433506044905009
$ws.Range("F19").Value = 0.01579925650557621
$ws.Range("H19").Value = 0.2235130111524164
$ws.Range("I19").Value = 0.07388475836431227
$ws.Range("J19").Value = 0.3559479553903346
$ws.Range("K19").Value = 0.1101301115241636
$ws.Range("M19").Value = 0.02184014869888476
$ws.Range("N19").Value = 0.0009293680297397769
$ws.Range("O19").Value = 0.06389405204460967
$ws.Range("S19").Value = 0.1340613382899628
